# Swap the presentation's applied theme from the "Integral" (Red Violet)
# color scheme back to the default "Office Theme" color scheme.
#
# (The OOXML-level change being reproduced renames/reswaps the two
# ppt/theme/themeN.xml parts so that the theme actually used by the
# slide master - and therefore by every slide - carries the "Office
# Theme" / "Office" colors that used to live in ppt/theme/theme1.xml,
# while the Red-Violet/"Integral" colors move to the other theme part.
# From the PowerPoint object model the reachable, user-visible half of
# that edit is the slide master's live theme color scheme, which this
# script updates to the Office Theme palette.)

function ConvertTo-BGR([string]$hex) {
    # MsoThemeColorSchemeIndex colors (and the legacy ColorScheme RGB
    # property) are surfaced as decimal 0x00BBGGRR, i.e. blue/green/red
    # byte order rather than the usual RRGGBB hex most people read.
    $v = [Convert]::ToInt32($hex, 16)
    $r = ($v -band 0xFF0000) -shr 16
    $g = ($v -band 0x00FF00) -shr 8
    $b = ($v -band 0x0000FF)
    return ($b -shl 16) -bor ($g -shl 8) -bor $r
}

$p = $ppt.ActivePresentation

# The 12-slot theme color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) in that exact order - matches a:clrScheme child order.
$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$tcs = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $tcs.Colors($i).RGB = ConvertTo-BGR $officeThemeColors[$i - 1]
}
